$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 29 (2025Q3) metrics: total_customers, returning_customers,
# new_customers, recurrence_rate
$ws.Range("C29").Value = 112
$ws.Range("D29").Value = 21
$ws.Range("E29").Value = 91
$ws.Range("F29").Value = 3.614457831325301
